$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to D, E, G columns for affected rows so that
# numeric-looking values (prices, percentages, hour) are stored as literal text,
# matching the original inlineStr string cells rather than being auto-converted
# to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "330.11"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.70%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "13"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.36"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.03%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "13"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.503"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.32%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "13"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08029"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.39%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "13"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.064"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "8.66%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "13"
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9543"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.79%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "13"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.556"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.10%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "13"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1140"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.22%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "13"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1904"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.85%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "13"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.18"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.91%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "13"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09902"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.59%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "13"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04883"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "11.75%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "13"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1064"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.34%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "13"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001263"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.62%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "13"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04086"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.26%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "13"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005892"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.13%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "13"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.391"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.53%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "13"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.406"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.38%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "13"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.99%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "13"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1382"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.46%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "13"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2580"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.69%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "13"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.36%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "13"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004355"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-3.44%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "13"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001199"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.91%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "13"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003741"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.29%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "13"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "13"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "13"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "13"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "13"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "13"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "13"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "13"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "13"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "13"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "13"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "13"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02596"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-1.85%"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "13"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05822"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.79%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "13"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007564"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.35%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "13"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1404"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.28%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "13"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007309"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "3.23%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "13"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002014"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.19%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "13"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008196"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.07%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "13"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00007024"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "1.47%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "13"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000749"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.16%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "13"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005799"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.21%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "13"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003528"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "55.26%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "13"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-8.72%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "13"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002099"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.16%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "13"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.16%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "13"
